$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add a new worksheet right after Sheet1 and name it "Sheet2"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Populate the new sheet
$ws2.Range("A1").Value = "validation base with gap to select features"
$ws2.Range("B1").Value = 2.03591244965226

# Match the column widths observed after the edit
$ws2.Columns.Item(1).ColumnWidth = 38.166666666666664
$ws2.Columns.Item(2).ColumnWidth = 12.7

# Portrait page orientation, like Sheet1
$ws2.PageSetup.Orientation = 1

# Sheet1's old selection moves to B11, while Sheet1 is still active
[void]$ws1.Range("B11").Select()

# Sheet2 becomes the active (selected) sheet/tab, with B1 selected
[void]$ws2.Activate()
[void]$ws2.Range("B1").Select()
